$wb = $excel.ActiveWorkbook

# ---- Sheet "Prix Spot": insert a new date column (30-nov) before column EC ----
$wsPrix = $wb.Worksheets.Item("Prix Spot")
$wsPrix.Range("EC1").EntireColumn.Insert()
$wsPrix.Range("EC1").Value2 = "30-nov"
$wsPrix.Range("EC2:EC25").Value2 = "-"

# ---- Sheet "Gaz": append a new row for 2025-11-28 ----
$wsGaz = $wb.Worksheets.Item("Gaz")
# Write the date as text (leading apostrophe forces text), then copy the
# cell format from the previous date cell so the new cell ends up with the
# same (default) style as every other date cell in the column.
$wsGaz.Range("A163").Value2 = "'2025-11-28"
$wsGaz.Range("A162").Copy()
$wsGaz.Range("A163").PasteSpecial(-4122)
$wsGaz.Application.CutCopyMode = $false
$wsGaz.Range("B163").Value2 = 27.5

# ---- Sheet "CO2": append a new row for 2025-11-28 (no price published yet) ----
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A163").Value2 = "'2025-11-28"
$wsCo2.Range("A162").Copy()
$wsCo2.Range("A163").PasteSpecial(-4122)
$wsCo2.Application.CutCopyMode = $false
